# Restored from revision of admin on 04/22/2021 11:55:09 AM.TEST Author: admin. Type: SAVE.
# Update the "Integer min" value for rule R20 (row 10) from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
